$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Fix typos / accents in the description column (A2:A14)
$ws.Range("A4").Value = "Mise en place d'Eedomus"

# Fix the header typo "Descritpion" -> "Description"
$ws.Range("A1").Value = "Description"

$ws.Range("A5").Value = "Réinstallation Raspberry pi et domoticz"
$ws.Range("A6").Value = "Théorie câblage plus début de la maquette"
$ws.Range("A7").Value = "Maquette presque terminée, circuit 12v restant"
$ws.Range("A8").Value = "Fin de la maquette sans qubino, début de l'ajout du qubino"
$ws.Range("A10").Value = "Ajout des appareils enocean et début des scénarios"
$ws.Range("A15").Value = "Finalisation de la documentation"

# Update the active selection cell to mirror the authored file
$ws.Range("I12").Select()
